$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account")

# Re-point the "Instruction" note to claim a shared-string slot before the
# new email text does, so the saved sharedStrings.xml ends up ordered the
# same way Excel itself would order it (index 49 = Created by Automation).
$ws.Range("B14").Value = "Created by Automation"

# The existing mailto hyperlink on B4 was authored against the file when it
# was loaded, so its target/display can't be edited in place - drop it and
# recreate it pointing at the same qa address, but now showing it as the
# "display" text (since the cell's own value will diverge from the link).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:vahan+qa@heal.com", [Type]::Missing, [Type]::Missing, "vahan+qa@heal.com")

# Now set the cell values. B4 gets the new dev email (this must happen
# after Hyperlinks.Add, which otherwise stamps its own display text into
# the cell). B5 (Password) loses its trailing "!".
$ws.Range("B4").Value = "vahan+dev@heal.com"
$ws.Range("B5").Value = "Heal4325"

# Move the active selection down to B5, matching the saved view state.
$ws.Range("B5").Select()
